$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 456, shifting existing rows 456..496 down to 457..497
$ws.Rows.Item(456).Insert()

# Populate the newly inserted row 456 with the new record's data
$ws.Cells.Item(456, 1).Value = 11
$ws.Cells.Item(456, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(456, 3).Value = "Bíobío"
$ws.Cells.Item(456, 4).Value = 44890
$ws.Cells.Item(456, 5).Value = 8
$ws.Cells.Item(456, 6).Value = 100114001
$ws.Cells.Item(456, 7).Value = "Papa"
$ws.Cells.Item(456, 8).Value = "Asterix"
$ws.Cells.Item(456, 9).Value = "1a (cosecha)"
$ws.Cells.Item(456, 10).Value = 200
$ws.Cells.Item(456, 11).Value = 14000
$ws.Cells.Item(456, 12).Value = 15000
$ws.Cells.Item(456, 13).Value = 14500
$ws.Cells.Item(456, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(456, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(456, 16).Value = 580
$ws.Cells.Item(456, 17).Value = 25
$ws.Cells.Item(456, 18).Value = "Hortaliza"
